$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analysis date (A2:A5) from 2025-12-01 to 2025-12-03
# Keep the cells as plain text (as they were originally) instead of letting
# Excel auto-convert the string into a date serial number.
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2:A5").Value = "2025-12-03"
$ws.Range("A2:A5").ClearFormats()

# Row 2 - IBM
$ws.Range("D2").Value = 304.73
$ws.Range("E2").Value = 41.9
$ws.Range("F2").Value = 0.2
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 59.8
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 65.32892478746797
$ws.Range("O2").Value = "🟢 상승 우위 (다소 완화)"

# Row 3 - QBTS
$ws.Range("D3").Value = 22.17
$ws.Range("E3").Value = 28.5
$ws.Range("F3").Value = -4.09
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 83
$ws.Range("I3").Value = 90
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 58.6
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 65.32892478746797
$ws.Range("O3").Value = "🟢 상승 우위 (다소 완화)"

# Row 4 - RGTI
$ws.Range("D4").Value = 23.93
$ws.Range("E4").Value = 29.9
$ws.Range("F4").Value = -9.94
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 66
$ws.Range("I4").Value = 83
$ws.Range("J4").Value = 86
$ws.Range("K4").Value = 55.8
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 65.32892478746797
$ws.Range("O4").Value = "🟢 상승 우위 (다소 완화)"

# Row 5 - IONQ
$ws.Range("D5").Value = 47.5
$ws.Range("E5").Value = 39.2
$ws.Range("F5").Value = 1.59
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 53
$ws.Range("J5").Value = 66
$ws.Range("K5").Value = 49.8
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 65.32892478746797
$ws.Range("O5").Value = "🟢 상승 우위 (다소 완화)"
